# The "связи между объектами" bullet list originally paired every object with
# itself (ВторойОбъект/ВторойОбъект, ТретийОбъект/ТретийОбъект, ...). The edit
# rewires each bullet to describe the relationship between two *different*
# objects instead, while keeping the relationship-kind wording unchanged.
#
# We locate each target paragraph with a short, unique substring (so the
# lookup is robust to any earlier formatting) and then replace the whole
# paragraph's text in one shot via Range.Text, which avoids Word's
# Find/Replace AutoCorrect machinery (it was silently turning the straight
# quotes into curly ones).

$d = $word.ActiveDocument
$paras = $d.Paragraphs
$count = $paras.Count

$markers = @(
    'объекта "ВторойОбъект" и объекта "ВторойОбъект"',
    'объекта "ТретийОбъект" и объекта "ТретийОбъект"',
    'объекта "Четвертый объект" и объекта "Четвертый объект"',
    'объекта "some object" и объекта "some object"'
)
$replacements = @(
    '— для связи объекта "ВторойОбъект" и объекта "ТретийОбъект" справедливо, что для одного объекта "ВторойОбъект" может существовать много объектов  "ТретийОбъект", и наоборот, для одного объекта "ТретийОбъект" может существовать только один объект   "ВторойОбъект", т.е. связь типа «один-ко-многим»;',
    '— для связи объекта "ТретийОбъект" и объекта "some object" справедливо, что для одного объекта "ТретийОбъект" может существовать много объектов  "some object", и наоборот, для одного объекта "some object" может существовать много объектов  "ТретийОбъект", т.е. связь типа «многие-ко-многим»;',
    '— для связи объекта "Четвертый объект" и объекта "some object" справедливо, что для одного объекта "some object" может существовать много объектов  "Четвертый объект", и наоборот, для одного объекта "Четвертый объект" может существовать только один объект   "some object", т.е. связь типа «один-ко-многим»;',
    '— для связи объекта "some object" и объекта "ВторойОбъект" справедливо, что для одного объекта "some object" может существовать только один объект   "ВторойОбъект", и наоборот, для одного объекта "ВторойОбъект" может существовать только один объект   "some object", т.е. связь типа «один-ко-одному»;'
)

for ($m = 0; $m -lt $markers.Length; $m++) {
    $marker = $markers[$m]
    $replacement = $replacements[$m]
    $matched = $false
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($marker)) {
            $p.Range.Text = $replacement
            $matched = $true
            break
        }
    }
    if (-not $matched) {
        throw "Could not find target paragraph for marker: $marker"
    }
    Write-Host "marker $m matched: $matched"
}
